$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix pre-existing cells ---
# Q61: 2 -> 0
$ws.Range("Q61").Value = 0
# Q69: 1 -> 0
$ws.Range("Q69").Value = 0
# O917: 0 -> 2
$ws.Range("O917").Value = 2
# R919 / R920: blank inlineStr -> numeric 0
$ws.Range("R919").Value = 0
$ws.Range("R920").Value = 0

# --- Append 20 new weekly rows (921-940) ---
# Column A uses the same custom date-time number format as the existing data rows.
$ws.Range("A921:A940").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 921
$ws.Range("A921").Value = 45474
$ws.Range("B921").Value = 6750.10009765625
$ws.Range("C921").Value = 6969
$ws.Range("D921").Value = 6642
$ws.Range("E921").Value = 6809.5
$ws.Range("F921").Value = 6798.130859375
$ws.Range("G921").Value = 1278020
$ws.Range("H921").Value = 2024
$ws.Range("I921").Value = 7
$ws.Range("J921").Value = 1
$ws.Range("K921").Value = 0
$ws.Range("L921").Value = 0
$ws.Range("M921").Value = 0
$ws.Range("N921").Value = 27
$ws.Range("O921").Value = 0
$ws.Range("P921").Value = 0
$ws.Range("Q921").Value = 0

# Row 922
$ws.Range("A922").Value = 45481
$ws.Range("B922").Value = 6874.5
$ws.Range("C922").Value = 7063
$ws.Range("D922").Value = 6569.10009765625
$ws.Range("E922").Value = 6861.5
$ws.Range("F922").Value = 6850.0439453125
$ws.Range("G922").Value = 1693031
$ws.Range("H922").Value = 2024
$ws.Range("I922").Value = 7
$ws.Range("J922").Value = 8
$ws.Range("K922").Value = 0
$ws.Range("L922").Value = 0
$ws.Range("M922").Value = 0
$ws.Range("N922").Value = 28
$ws.Range("O922").Value = 0
$ws.Range("P922").Value = 0
$ws.Range("Q922").Value = 0

# Row 923
$ws.Range("A923").Value = 45488
$ws.Range("B923").Value = 6870
$ws.Range("C923").Value = 7069.9501953125
$ws.Range("D923").Value = 6741.25
$ws.Range("E923").Value = 6827.60009765625
$ws.Range("F923").Value = 6816.20068359375
$ws.Range("G923").Value = 829076
$ws.Range("H923").Value = 2024
$ws.Range("I923").Value = 7
$ws.Range("J923").Value = 15
$ws.Range("K923").Value = 0
$ws.Range("L923").Value = 0
$ws.Range("M923").Value = 0
$ws.Range("N923").Value = 29
$ws.Range("O923").Value = 0
$ws.Range("P923").Value = 0
$ws.Range("Q923").Value = 1

# Row 924
$ws.Range("A924").Value = 45495
$ws.Range("B924").Value = 6824.9501953125
$ws.Range("C924").Value = 7210
$ws.Range("D924").Value = 6706.7998046875
$ws.Range("E924").Value = 7187.39990234375
$ws.Range("F924").Value = 7175.39990234375
$ws.Range("G924").Value = 1211226
$ws.Range("H924").Value = 2024
$ws.Range("I924").Value = 7
$ws.Range("J924").Value = 22
$ws.Range("K924").Value = 0
$ws.Range("L924").Value = 0
$ws.Range("M924").Value = 0
$ws.Range("N924").Value = 30
$ws.Range("O924").Value = 0
$ws.Range("P924").Value = 0
$ws.Range("Q924").Value = 0

# Row 925
$ws.Range("A925").Value = 45502
$ws.Range("B925").Value = 7244.9501953125
$ws.Range("C925").Value = 7327.75
$ws.Range("D925").Value = 6844.7998046875
$ws.Range("E925").Value = 7220.5
$ws.Range("F925").Value = 7208.44482421875
$ws.Range("G925").Value = 2220538
$ws.Range("H925").Value = 2024
$ws.Range("I925").Value = 7
$ws.Range("J925").Value = 29
$ws.Range("K925").Value = 0
$ws.Range("L925").Value = 0
$ws.Range("M925").Value = 0
$ws.Range("N925").Value = 31
$ws.Range("O925").Value = 0
$ws.Range("P925").Value = 0
$ws.Range("Q925").Value = 0

# Row 926
$ws.Range("A926").Value = 45509
$ws.Range("B926").Value = 7000.0498046875
$ws.Range("C926").Value = 7250
$ws.Range("D926").Value = 6782
$ws.Range("E926").Value = 7214.5
$ws.Range("F926").Value = 7214.5
$ws.Range("G926").Value = 1615923
$ws.Range("H926").Value = 2024
$ws.Range("I926").Value = 8
$ws.Range("J926").Value = 5
$ws.Range("K926").Value = 0
$ws.Range("L926").Value = 0
$ws.Range("M926").Value = 0
$ws.Range("N926").Value = 32
$ws.Range("O926").Value = 0
$ws.Range("P926").Value = 0
$ws.Range("Q926").Value = 0

# Row 927
$ws.Range("A927").Value = 45516
$ws.Range("B927").Value = 7300
$ws.Range("C927").Value = 7440
$ws.Range("D927").Value = 7039.0498046875
$ws.Range("E927").Value = 7417.85009765625
$ws.Range("F927").Value = 7417.85009765625
$ws.Range("G927").Value = 1150359
$ws.Range("H927").Value = 2024
$ws.Range("I927").Value = 8
$ws.Range("J927").Value = 12
$ws.Range("K927").Value = 0
$ws.Range("L927").Value = 0
$ws.Range("M927").Value = 0
$ws.Range("N927").Value = 33
$ws.Range("O927").Value = 0
$ws.Range("P927").Value = 0
$ws.Range("Q927").Value = 0

# Row 928
$ws.Range("A928").Value = 45523
$ws.Range("B928").Value = 7417.85009765625
$ws.Range("C928").Value = 7570.0498046875
$ws.Range("D928").Value = 7314.0498046875
$ws.Range("E928").Value = 7392.60009765625
$ws.Range("F928").Value = 7392.60009765625
$ws.Range("G928").Value = 826152
$ws.Range("H928").Value = 2024
$ws.Range("I928").Value = 8
$ws.Range("J928").Value = 19
$ws.Range("K928").Value = 0
$ws.Range("L928").Value = 0
$ws.Range("M928").Value = 0
$ws.Range("N928").Value = 34
$ws.Range("O928").Value = 0
$ws.Range("P928").Value = 0
$ws.Range("Q928").Value = 0

# Row 929
$ws.Range("A929").Value = 45530
$ws.Range("B929").Value = 7387.64990234375
$ws.Range("C929").Value = 7725
$ws.Range("D929").Value = 7387.64990234375
$ws.Range("E929").Value = 7678.7998046875
$ws.Range("F929").Value = 7678.7998046875
$ws.Range("G929").Value = 1098543
$ws.Range("H929").Value = 2024
$ws.Range("I929").Value = 8
$ws.Range("J929").Value = 26
$ws.Range("K929").Value = 0
$ws.Range("L929").Value = 0
$ws.Range("M929").Value = 0
$ws.Range("N929").Value = 35
$ws.Range("O929").Value = 0
$ws.Range("P929").Value = 0
$ws.Range("Q929").Value = 0

# Row 930
$ws.Range("A930").Value = 45537
$ws.Range("B930").Value = 7679.9501953125
$ws.Range("C930").Value = 7746.5
$ws.Range("D930").Value = 7264.0498046875
$ws.Range("E930").Value = 7424.89990234375
$ws.Range("F930").Value = 7424.89990234375
$ws.Range("G930").Value = 809628
$ws.Range("H930").Value = 2024
$ws.Range("I930").Value = 9
$ws.Range("J930").Value = 2
$ws.Range("K930").Value = 0
$ws.Range("L930").Value = 0
$ws.Range("M930").Value = 0
$ws.Range("N930").Value = 36
$ws.Range("O930").Value = 0
$ws.Range("P930").Value = 0
$ws.Range("Q930").Value = 0

# Row 931
$ws.Range("A931").Value = 45544
$ws.Range("B931").Value = 7331
$ws.Range("C931").Value = 7830
$ws.Range("D931").Value = 7296.4501953125
$ws.Range("E931").Value = 7729.0498046875
$ws.Range("F931").Value = 7729.0498046875
$ws.Range("G931").Value = 987360
$ws.Range("H931").Value = 2024
$ws.Range("I931").Value = 9
$ws.Range("J931").Value = 9
$ws.Range("K931").Value = 0
$ws.Range("L931").Value = 0
$ws.Range("M931").Value = 0
$ws.Range("N931").Value = 37
$ws.Range("O931").Value = 0
$ws.Range("P931").Value = 0
$ws.Range("Q931").Value = 0

# Row 932
$ws.Range("A932").Value = 45551
$ws.Range("B932").Value = 7729.0498046875
$ws.Range("C932").Value = 8151
$ws.Range("D932").Value = 7570.0498046875
$ws.Range("E932").Value = 8106.2001953125
$ws.Range("F932").Value = 8106.2001953125
$ws.Range("G932").Value = 1801547
$ws.Range("H932").Value = 2024
$ws.Range("I932").Value = 9
$ws.Range("J932").Value = 16
$ws.Range("K932").Value = 0
$ws.Range("L932").Value = 0
$ws.Range("M932").Value = 0
$ws.Range("N932").Value = 38
$ws.Range("O932").Value = 0
$ws.Range("P932").Value = 0
$ws.Range("Q932").Value = 0

# Row 933
$ws.Range("A933").Value = 45558
$ws.Range("B933").Value = 8111.5
$ws.Range("C933").Value = 8260.349609375
$ws.Range("D933").Value = 7834.10009765625
$ws.Range("E933").Value = 8170.75
$ws.Range("F933").Value = 8170.75
$ws.Range("G933").Value = 1223261
$ws.Range("H933").Value = 2024
$ws.Range("I933").Value = 9
$ws.Range("J933").Value = 23
$ws.Range("K933").Value = 0
$ws.Range("L933").Value = 0
$ws.Range("M933").Value = 0
$ws.Range("N933").Value = 39
$ws.Range("O933").Value = 0
$ws.Range("P933").Value = 0
$ws.Range("Q933").Value = 0

# Row 934
$ws.Range("A934").Value = 45565
$ws.Range("B934").Value = 8113.7001953125
$ws.Range("C934").Value = 8308.900390625
$ws.Range("D934").Value = 8010.64990234375
$ws.Range("E934").Value = 8198.650390625
$ws.Range("F934").Value = 8198.650390625
$ws.Range("G934").Value = 907672
$ws.Range("H934").Value = 2024
$ws.Range("I934").Value = 9
$ws.Range("J934").Value = 30
$ws.Range("K934").Value = 0
$ws.Range("L934").Value = 0
$ws.Range("M934").Value = 0
$ws.Range("N934").Value = 40
$ws.Range("O934").Value = 0
$ws.Range("P934").Value = 0
$ws.Range("Q934").Value = 0

# Row 935
$ws.Range("A935").Value = 45572
$ws.Range("B935").Value = 8308.900390625
$ws.Range("C935").Value = 8472
$ws.Range("D935").Value = 7910.0498046875
$ws.Range("E935").Value = 8363.2998046875
$ws.Range("F935").Value = 8363.2998046875
$ws.Range("G935").Value = 1032010
$ws.Range("H935").Value = 2024
$ws.Range("I935").Value = 10
$ws.Range("J935").Value = 7
$ws.Range("K935").Value = 0
$ws.Range("L935").Value = 0
$ws.Range("M935").Value = 0
$ws.Range("N935").Value = 41
$ws.Range("O935").Value = 1
$ws.Range("P935").Value = 0
$ws.Range("Q935").Value = 0

# Row 936
$ws.Range("A936").Value = 45579
$ws.Range("B936").Value = 8389.9501953125
$ws.Range("C936").Value = 8399.25
$ws.Range("D936").Value = 7910.7001953125
$ws.Range("E936").Value = 8088.85009765625
$ws.Range("F936").Value = 8088.85009765625
$ws.Range("G936").Value = 843011
$ws.Range("H936").Value = 2024
$ws.Range("I936").Value = 10
$ws.Range("J936").Value = 14
$ws.Range("K936").Value = 0
$ws.Range("L936").Value = 0
$ws.Range("M936").Value = 0
$ws.Range("N936").Value = 42
$ws.Range("O936").Value = 0
$ws.Range("P936").Value = 0
$ws.Range("Q936").Value = 0

# Row 937
$ws.Range("A937").Value = 45586
$ws.Range("B937").Value = 8117
$ws.Range("C937").Value = 8159.75
$ws.Range("D937").Value = 7434.0498046875
$ws.Range("E937").Value = 7597.64990234375
$ws.Range("F937").Value = 7597.64990234375
$ws.Range("G937").Value = 1086589
$ws.Range("H937").Value = 2024
$ws.Range("I937").Value = 10
$ws.Range("J937").Value = 21
$ws.Range("K937").Value = 0
$ws.Range("L937").Value = 0
$ws.Range("M937").Value = 0
$ws.Range("N937").Value = 43
$ws.Range("O937").Value = 0
$ws.Range("P937").Value = 0
$ws.Range("Q937").Value = 0

# Row 938
$ws.Range("A938").Value = 45593
$ws.Range("B938").Value = 7600.0498046875
$ws.Range("C938").Value = 7856.4501953125
$ws.Range("D938").Value = 7406.5498046875
$ws.Range("E938").Value = 7468.4501953125
$ws.Range("F938").Value = 7468.4501953125
$ws.Range("G938").Value = 778301
$ws.Range("H938").Value = 2024
$ws.Range("I938").Value = 10
$ws.Range("J938").Value = 28
$ws.Range("K938").Value = 0
$ws.Range("L938").Value = 0
$ws.Range("M938").Value = 0
$ws.Range("N938").Value = 44
$ws.Range("O938").Value = 0
$ws.Range("P938").Value = 0
$ws.Range("Q938").Value = 0

# Row 939
$ws.Range("A939").Value = 45600
$ws.Range("B939").Value = 7450
$ws.Range("C939").Value = 8045.85009765625
$ws.Range("D939").Value = 7370.5498046875
$ws.Range("E939").Value = 7654.60009765625
$ws.Range("F939").Value = 7654.60009765625
$ws.Range("G939").Value = 1616311
$ws.Range("H939").Value = 2024
$ws.Range("I939").Value = 11
$ws.Range("J939").Value = 4
$ws.Range("K939").Value = 0
$ws.Range("L939").Value = 0
$ws.Range("M939").Value = 0
$ws.Range("N939").Value = 45
$ws.Range("O939").Value = 0
$ws.Range("P939").Value = 0
$ws.Range("Q939").Value = 0

# Row 940
$ws.Range("A940").Value = 45607
$ws.Range("B940").Value = 7670
$ws.Range("C940").Value = 8100
$ws.Range("D940").Value = 7550
$ws.Range("E940").Value = 7768.2001953125
$ws.Range("F940").Value = 7768.2001953125
$ws.Range("G940").Value = 1273066
$ws.Range("H940").Value = 2024
$ws.Range("I940").Value = 11
$ws.Range("J940").Value = 11
$ws.Range("K940").Value = 0
$ws.Range("L940").Value = 0
$ws.Range("M940").Value = 0
$ws.Range("N940").Value = 46
$ws.Range("O940").Value = 0
$ws.Range("P940").Value = 0
$ws.Range("Q940").Value = 0

